$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new cells in the same order the shared strings were interned
# (codigo row66, iniciativa row67, codigo row67, iniciativa row66)
$ws.Range("A66").Value = "9468204"
$ws.Range("B67").Value = "Elaboração de projetos rodoviários - Brumadinho-Mário Campos-BR381"
$ws.Range("A67").Value = "9468058"
$ws.Range("B66").Value = "Melhoria da infraestrutura dos municípios - Fortalecimento do saneamento básico de Mário Campos"

$ws.Range("C66").Value = "IV"
$ws.Range("D66").Value = 50000000
$ws.Range("C67").Value = "IV"
$ws.Range("D67").Value = 10000000

# Row height for the new rows
$ws.Rows.Item(66).RowHeight = 30
$ws.Rows.Item(67).RowHeight = 30

# Update the view: scroll position and selection to match the new state
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("B67").Select()
